$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Skill Description") before the current SFIA Level column.
$ws.Columns.Item(2).Insert()

# Header row
$ws.Cells.Item(1, 1).Value = "SkillCode"
$ws.Cells.Item(1, 2).Value = "Skill Description"
$ws.Cells.Item(1, 3).Value = "SFIA Level"
$ws.Cells.Item(1, 4).Value = "Keycode"
$ws.Cells.Item(1, 5).Value = "Description"

# Map of SkillCode -> full Skill Description name
$names = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "SORC"       = "Sourcing"
    "SUPP"       = "Supplier management"
    "ITCM"       = "Contract management"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($names.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $names[$code]
    }
}
